$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text.Replace("✅ 1000 Bs = 8.94 = 36800.71 pesos", "✅ 1000 Bs = 9.12 = 37667.06 pesos")
$text = $text.Replace("✅ 36800.71 pesos = 8.88 = 958.63 Bs", "✅ 37667.06 pesos = 9.07 = 959.4 Bs")
$cellA1.Value = $text

# --- Update "tasas" sheet N10/O10/N12/O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 109.6
$wsTasas.Range("O10").Value = 4128.31
$wsTasas.Range("N12").Value = 4154
$wsTasas.Range("O12").Value = 105.805
